$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 0.5327
$ws.Range("C3").Value = 0.6011
$ws.Range("C4").Value = 0.2799
$ws.Range("C5").Value = 0.4121
$ws.Range("C6").Value = -0.1786
